# Structure revamp: normalize the "document" (A, F) columns to lowercase
# headline text while giving each row's "h1" (C) cell its own
# Title-Case display string (rows 3-8 already followed this pattern;
# row 2 is brought in line with it here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "saas product design principle"
$ws.Range("F2").Value = "saas product design principle"
$ws.Range("C2").Value = "Saas Product Design Principle"

$ws.Range("A3").Value = "intelligent b-side product analysis and design ideas"
$ws.Range("F3").Value = "intelligent b-side product analysis and design ideas"

$ws.Range("A4").Value = "how to design a corporate system ui"
$ws.Range("F4").Value = "how to design a corporate system ui"

$ws.Range("A5").Value = "how to design a 2b data table"
$ws.Range("F5").Value = "how to design a 2b data table"

$ws.Range("A6").Value = "how to break down complex process"
$ws.Range("F6").Value = "how to break down complex process"

$ws.Range("A7").Value = "business requirement management for corporate systems"
$ws.Range("F7").Value = "business requirement management for corporate systems"
